$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New issue row (row 47): Description, Status, Status Date
$ws.Range("B47").Value = "Allow GUI to select the calibration file on the PNA"
$ws.Range("C47").Value = "Open"

# Copy the existing date formatting from D46 so D47 reuses the same style
# (numFmtId 14, "m/d/yyyy") instead of Excel creating a brand new style.
$ws.Range("D46").Copy()
$ws.Range("D47").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D47").Value = 43451

# Update the view state to match where the user ended up after the edit.
$ws.Application.ActiveWindow.ScrollRow = 27
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("B48").Select()
